# Update "horarios 141" workbook with the latest scrape results (01:53:21).

$wb = $excel.ActiveWorkbook

$newTime = "01:53:21"

# --- Sheet 1: LP1912 -------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: $newTime"
$ws1.Range("A3").Value = "Total filas: 3"

# Existing row 6 gets updated values (now the "15_ABASTO" @ 03:02 entry).
$ws1.Cells.Item(6, 1).Value = $newTime
$ws1.Cells.Item(6, 2).Value = "03:02"
$ws1.Cells.Item(6, 3).Value = "15_ABASTO"
$ws1.Cells.Item(6, 4).Value = 69
$ws1.Cells.Item(6, 5).Value = "LP1912"

# Existing row 7 gets updated values (still "15_ABASTO", now @ 03:02).
$ws1.Cells.Item(7, 1).Value = $newTime
$ws1.Cells.Item(7, 2).Value = "03:02"
$ws1.Cells.Item(7, 3).Value = "15_ABASTO"
$ws1.Cells.Item(7, 4).Value = 69
$ws1.Cells.Item(7, 5).Value = "LP1912"

# New row 8: additional arrival added during this scrape.
$ws1.Cells.Item(8, 1).Value = $newTime
$ws1.Cells.Item(8, 2).Value = "03:48"
$ws1.Cells.Item(8, 3).Value = "14_ABASTO"
$ws1.Cells.Item(8, 4).Value = 115
$ws1.Cells.Item(8, 5).Value = "LP1912"

# --- Sheet 2: LP1912-215 ----------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: $newTime"

# --- Sheet 3: 6203-6173 ------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: $newTime"
